# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever both appear together as "System, dnasr281@gmail.com".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
